$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 6593.25  # was 7999.6
$ws.Range("I100").Value = 5949.75  # was 7666.3335
$ws.Range("J100").Value = 7236.75  # was 8499.5
$ws.Range("K100").Value = 5949.75  # was 7666.3335
$ws.Range("L100").Value = 7236.75  # was 8499.5
$ws.Range("M100").Value = -5408.75  # was -7125.3335
$ws.Range("N100").Value = -8318.75  # was -9581.5
$ws.Range("H107").Value = 787.36365  # was 770.44446
$ws.Range("I107").Value = 751  # was 718.8570999999999
$ws.Range("K107").Value = 751  # was 718.8570999999999
$ws.Range("M107").Value = 1169  # was 1201.1429
$ws.Range("H112").Value = 1740.1154  # was 1874.3334
$ws.Range("I112").Value = 1399.6  # was 1365.5
$ws.Range("J112").Value = 1821.1904  # was 2128.75
$ws.Range("K112").Value = 4198.799999999999  # was 4096.5
$ws.Range("L112").Value = 5463.5712  # was 6386.25
$ws.Range("M112").Value = -3090.799999999999  # was -2988.5
$ws.Range("N112").Value = -7679.5712  # was -8602.25
$ws.Range("H113").Value = 4470.222  # was 4367.727
$ws.Range("I113").Value = 4133  # was 4086.6667
$ws.Range("J113").Value = 4891.75  # was 4705
$ws.Range("K113").Value = 4133  # was 4086.6667
$ws.Range("L113").Value = 4891.75  # was 4705
$ws.Range("M113").Value = -879  # was -832.6667000000002
$ws.Range("N113").Value = -11399.75  # was -11213
$ws.Range("H129").Value = 2353.6897  # was 2353
$ws.Range("I129").Value = 874.53845  # was 826.3570999999999
$ws.Range("J129").Value = 3555.5  # was 3777.8667
$ws.Range("K129").Value = 2623.61535  # was 2479.0713
$ws.Range("L129").Value = 10666.5  # was 11333.6001
$ws.Range("M129").Value = 2376.38465  # was 2520.9287
$ws.Range("N129").Value = -20666.5  # was -21333.6001
$ws.Range("H132").Value = 3035.1765  # was 3210.5
$ws.Range("I132").Value = 2850.5  # was 3025.2
$ws.Range("K132").Value = 8551.5  # was 9075.599999999999
$ws.Range("M132").Value = -6021.5  # was -6545.599999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1816.6923  # was 1989.8334
$ws.Range("I2").Value = 1801.4166  # was 2007.8
$ws.Range("J2").Value = 2000  # was 1900
$ws.Range("K2").Value = 1801.4166  # was 2007.8
$ws.Range("L2").Value = 2000  # was 1900
$ws.Range("M2").Value = -1688.4166  # was -1894.8
$ws.Range("N2").Value = -2226  # was -2126
$ws.Range("H27").Value = 2500  # was 1833.3334
$ws.Range("J27").Value = 2500  # was 1833.3334
$ws.Range("L27").Value = 2500  # was 1833.3334
$ws.Range("N27").Value = -2868  # was -2201.3334
$ws.Range("H61").Value = 1640  # was 1725
$ws.Range("I61").Value = 1640  # was 1725
$ws.Range("K61").Value = 1640  # was 1725
$ws.Range("M61").Value = -1428  # was -1513
$ws.Range("H74").Value = 1278.7  # was 1003.5714
$ws.Range("I74").Value = 1278.7  # was 1003.5714
$ws.Range("K74").Value = 1278.7  # was 1003.5714
$ws.Range("M74").Value = -404.7  # was -129.5714
$ws.Range("H77").Value = 1278.7  # was 1003.5714
$ws.Range("I77").Value = 1278.7  # was 1003.5714
$ws.Range("K77").Value = 6393.5  # was 5017.857
$ws.Range("M77").Value = -2025.5  # was -649.857
$ws.Range("H110").Value = 1415.7778  # was 2984.1875
$ws.Range("I110").Value = 391.7143  # was 392.42856
$ws.Range("K110").Value = 391.7143  # was 392.42856
$ws.Range("M110").Value = 1653.2857  # was 1652.57144
$ws.Range("H116").Value = 1816.6923  # was 1989.8334
$ws.Range("I116").Value = 1801.4166  # was 2007.8
$ws.Range("J116").Value = 2000  # was 1900
$ws.Range("K116").Value = 1801.4166  # was 2007.8
$ws.Range("L116").Value = 2000  # was 1900
$ws.Range("M116").Value = 492.5834  # was 286.2
$ws.Range("N116").Value = -6588  # was -6488
$ws.Range("H122").Value = 2998.5  # was 2997.6667
$ws.Range("I122").Value = 2998.5  # was 2997.6667
$ws.Range("K122").Value = 8995.5  # was 8993.000100000001
$ws.Range("M122").Value = -6545.5  # was -6543.000100000001
$ws.Range("H136").Value = 1640  # was 1725
$ws.Range("I136").Value = 1640  # was 1725
$ws.Range("K136").Value = 4920  # was 5175
$ws.Range("M136").Value = -2370  # was -2625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1816.6923  # was 1989.8334
$ws.Range("I3").Value = 1801.4166  # was 2007.8
$ws.Range("J3").Value = 2000  # was 1900
$ws.Range("K3").Value = 1801.4166  # was 2007.8
$ws.Range("L3").Value = 2000  # was 1900
$ws.Range("M3").Value = -1687.4166  # was -1893.8
$ws.Range("N3").Value = -2228  # was -2128
$ws.Range("H94").Value = 3332.3333  # was 3999
$ws.Range("I94").Value = 3332.3333  # was 3999
$ws.Range("K94").Value = 3332.3333  # was 3999
$ws.Range("M94").Value = -2881.3333  # was -3548

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1878.3  # was 1897.3
$ws.Range("J4").Value = 1733.2667  # was 1771.2667
$ws.Range("L4").Value = 5199.800099999999  # was 5313.800099999999
$ws.Range("N4").Value = -5423.800099999999  # was -5537.800099999999
$ws.Range("H118").Value = 1725  # was 1875
$ws.Range("I118").Value = 1750  # was 1875
$ws.Range("J118").Value = 1700  # was 0
$ws.Range("K118").Value = 5250  # was 5625
$ws.Range("L118").Value = 5100  # was 0
$ws.Range("M118").Value = -4007  # was -4382
$ws.Range("N118").Value = -7586  # was None
$ws.Range("H120").Value = 13500  # was 7975
$ws.Range("I120").Value = 10000  # was 2956.25
$ws.Range("K120").Value = 30000  # was 8868.75
$ws.Range("M120").Value = -25162  # was -4030.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 29855.572  # was 30000
$ws.Range("I62").Value = 29495  # was 0
$ws.Range("J62").Value = 29999.8  # was 30000
$ws.Range("K62").Value = 29495  # was 0
$ws.Range("L62").Value = 29999.8  # was 30000
$ws.Range("M62").Value = -28809  # was None
$ws.Range("N62").Value = -31371.8  # was -31372
$ws.Range("H65").Value = 29855.572  # was 30000
$ws.Range("I65").Value = 29495  # was 0
$ws.Range("J65").Value = 29999.8  # was 30000
$ws.Range("K65").Value = 88485  # was 0
$ws.Range("L65").Value = 89999.39999999999  # was 90000
$ws.Range("M65").Value = -85053  # was None
$ws.Range("N65").Value = -96863.39999999999  # was -96864
$ws.Range("H122").Value = 8336213  # was 8931525
$ws.Range("I122").Value = 9618092  # was 10419446
$ws.Range("K122").Value = 28854276  # was 31258338
$ws.Range("M122").Value = -28851826  # was -31255888
$ws.Range("H126").Value = 3874  # was 3998.6667
$ws.Range("I126").Value = 3498.6667  # was 3998.6667
$ws.Range("J126").Value = 5000  # was 0
$ws.Range("K126").Value = 10496.0001  # was 11996.0001
$ws.Range("L126").Value = 15000  # was 0
$ws.Range("M126").Value = -8026.000100000001  # was -9526.000100000001
$ws.Range("N126").Value = -19940  # was None

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3286.2856  # was 3125.5
$ws.Range("I7").Value = 2251  # was 2200.8
$ws.Range("K7").Value = 2251  # was 2200.8
$ws.Range("M7").Value = -2139  # was -2088.8
$ws.Range("H126").Value = 3286.2856  # was 3125.5
$ws.Range("I126").Value = 2251  # was 2200.8
$ws.Range("K126").Value = 6753  # was 6602.400000000001
$ws.Range("M126").Value = -4283  # was -4132.400000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 928300  # was 87885.71000000001
$ws.Range("I4").Value = 1392000  # was 153250
$ws.Range("J4").Value = 900  # was 733.3333
$ws.Range("K4").Value = 1392000  # was 153250
$ws.Range("L4").Value = 900  # was 733.3333
$ws.Range("M4").Value = -1391887  # was -153137
$ws.Range("N4").Value = -1126  # was -959.3333
$ws.Range("H96").Value = 2000  # was 1037.5
$ws.Range("I96").Value = 0  # was 716.6667
$ws.Range("K96").Value = 0  # was 716.6667
$ws.Range("M96").ClearContents()  # was 656.3333
$ws.Range("H107").Value = 965.625  # was 963.6667
$ws.Range("I107").Value = 1283  # was 1073.6666
$ws.Range("J107").Value = 775.2  # was 853.6667
$ws.Range("K107").Value = 3849  # was 3220.9998
$ws.Range("L107").Value = 2325.6  # was 2561.0001
$ws.Range("M107").Value = -1929  # was -1300.9998
$ws.Range("N107").Value = -6165.6  # was -6401.0001
$ws.Range("H112").Value = 50000  # was 49444
$ws.Range("J112").Value = 50000  # was 49444
$ws.Range("L112").Value = 50000  # was 49444
$ws.Range("N112").Value = -52954  # was -52398
$ws.Range("H126").Value = 2695.3125  # was 2908.9285
$ws.Range("I126").Value = 2394.2307  # was 2611.3635
$ws.Range("K126").Value = 7182.6921  # was 7834.0905
$ws.Range("M126").Value = -4712.6921  # was -5364.0905
$ws.Range("H136").Value = 2538  # was 2758.4285
$ws.Range("I136").Value = 2538  # was 2758.4285
$ws.Range("K136").Value = 7614  # was 8275.2855
$ws.Range("M136").Value = -5064  # was -5725.2855
